$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and row 47-50 reordering) per commit

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.180.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.46%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.829.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.20%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.00"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6155"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.46%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07342"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2895"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.15%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.19"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07638"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.43%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.829.30"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.07%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.978"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.42%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6710"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.22%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.44"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.52%  "

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.00%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.848"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.01%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.175.86"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.077.94"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.05%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.27"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.51%  "

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.47%  "

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.371"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +2.35%  "

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.05%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.76"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.31%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.521"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.38%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.31%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.63"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.29%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.489"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.63%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05817"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.49%  "

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.39%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.03%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.088"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.40%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.851"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.54%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.49%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7201"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.63%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.617"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.863"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.19%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.227.13"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.54%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01761"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8994"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.19%  "

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.10%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.996.52"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.02%  "

# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.14%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.49"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.08%  "

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000119"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.40%  "

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5042"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.82%  "

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.186"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.78%  "

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4035"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.81%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1162"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.27%  "

